# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Update the "Periodo Mora" values (column E, rows 16-21) of the account
# statement table: the old period list (2507,2506,2505,2504,2503,2502) is
# replaced with the new period list (2503,2504,2505,2506,2507,2508).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2503"
$ws.Range("E17").Value = "2504"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2506"
$ws.Range("E20").Value = "2507"
$ws.Range("E21").Value = "2508"
